# "Generate Report for Archive"
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F2, zh-cn!C2,
#   de-de!C2 all share the same string).
# - Narrow the status column(s) on all three sheets to their new width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status values (shared text) on every sheet that references it.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Resize the status columns to match the new (narrower) width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
